# The Dragon Boat Festival
#
# 1. Paragraph 1 ("今天星期六，下雨又是开心的一天") becomes justified.
# 2. A new paragraph carrying the old paragraph-2 text
#    ("今天星期天，晴天又是开心的一天") is inserted right after paragraph 1.
# 3. The (now) third / final paragraph - the one that still owns the
#    trailing _GoBack bookmark - gets its text swapped to
#    "今天星期一又要去上课了，难受".

$d = $word.ActiveDocument

# --- Step 2 first: duplicate paragraph 2's wording into a fresh paragraph
# inserted after paragraph 1. Doing this before the alignment tweak keeps
# the new paragraph from inheriting a w:jc it should not have.
$firstPara = $d.Paragraphs(1)
$firstPara.Range.InsertParagraphAfter()
$newPara = $d.Paragraphs(2)
$newPara.Range.Text = "今天星期天，晴天又是开心的一天"

# --- Step 1: justify the first paragraph only.
$firstPara.Format.Alignment = 3

# --- Step 3: change the wording of the final paragraph (still holding the
# bookmark) from the Sunday sentence to the Monday sentence. Scope the
# Find/Replace to that paragraph's own Range so the identical text that now
# also lives in paragraph 2 is left untouched.
$lastPara = $d.Paragraphs($d.Paragraphs.Count)
$lastPara.Range.Find.Execute("今天星期天，晴天又是开心的一天", $true, $false, $false, $false, $false, `
                              $true, 1, $false, "今天星期一又要去上课了，难受", 2) | Out-Null
